$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# C2: "trainingaudio/14_pokoto1.wav" -> "train2P2"
$ws.Range("C2").Value = "train2P2"

# C3: "trainingaudio/13_kopopi1.wav" -> "train2P2"
$ws.Range("C3").Value = "train2P2"
